# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 930.7692
$ws.Range("J17").Value = 1200
$ws.Range("L17").Value = 3600
$ws.Range("N17").Value = -3936
$ws.Range("H116").Value = 23819774
$ws.Range("I116").Value = 23812860
$ws.Range("J116").Value = 23822738
$ws.Range("K116").Value = 23812860
$ws.Range("L116").Value = 23822738
$ws.Range("M116").Value = -23809418
$ws.Range("N116").Value = -23829622
$ws.Range("H118").Value = 2566.923
$ws.Range("I118").Value = 1796.25
$ws.Range("J118").Value = 3800
$ws.Range("K118").Value = 5388.75
$ws.Range("L118").Value = 11400
$ws.Range("M118").Value = -3731.75
$ws.Range("N118").Value = -14714
$ws.Range("H132").Value = 5924347
$ws.Range("I132").Value = 1256543
$ws.Range("J132").Value = 37043040
$ws.Range("K132").Value = 3769629
$ws.Range("L132").Value = 111129120
$ws.Range("M132").Value = -3767099
$ws.Range("N132").Value = -111134180
$ws.Range("H137").Value = 12721442
$ws.Range("I137").Value = 3788809.2
$ws.Range("J137").Value = 24059016
$ws.Range("K137").Value = 11366427.6
$ws.Range("L137").Value = 72177048
$ws.Range("M137").Value = -11363877.6
$ws.Range("N137").Value = -72182148
$ws.Range("H138").Value = 4299.5977
$ws.Range("I138").Value = 4968.1875
$ws.Range("J138").Value = 4158.8423
$ws.Range("K138").Value = 14904.5625
$ws.Range("L138").Value = 12476.5269
$ws.Range("M138").Value = -9764.5625
$ws.Range("N138").Value = -22756.5269

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 255.55556
$ws.Range("I5").Value = 255.55556
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 255.55556
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -143.55556
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 4333.3335
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -25900
$ws.Range("H132").Value = 19774358
$ws.Range("I132").Value = 22708568
$ws.Range("J132").Value = 9617475
$ws.Range("K132").Value = 68125704
$ws.Range("L132").Value = 28852425
$ws.Range("M132").Value = -68123174
$ws.Range("N132").Value = -28857485

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 255.55556
$ws.Range("I4").Value = 255.55556
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 255.55556
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -140.55556
$ws.Range("N4").ClearContents()
$ws.Range("H99").Value = 2126.6667
$ws.Range("I99").Value = 1216.3572
$ws.Range("J99").Value = 2636.44
$ws.Range("K99").Value = 1216.3572
$ws.Range("L99").Value = 2636.44
$ws.Range("M99").Value = 281.6428000000001
$ws.Range("N99").Value = -5632.440000000001
$ws.Range("H107").Value = 1430388.9
$ws.Range("I107").Value = 1430388.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1430388.9
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1428468.9
$ws.Range("N107").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 976.9231
$ws.Range("I22").Value = 166.77777
$ws.Range("K22").Value = 166.77777
$ws.Range("M22").Value = 183.22223
$ws.Range("H50").Value = 17758
$ws.Range("J50").Value = 17758
$ws.Range("L50").Value = 17758
$ws.Range("N50").Value = -19008
$ws.Range("H51").Value = 17433
$ws.Range("J51").Value = 17433
$ws.Range("L51").Value = 17433
$ws.Range("N51").Value = -18905
$ws.Range("H58").Value = 1757048.5
$ws.Range("I58").Value = 14074.333
$ws.Range("J58").Value = 4133831.5
$ws.Range("K58").Value = 14074.333
$ws.Range("L58").Value = 4133831.5
$ws.Range("M58").Value = -13871.333
$ws.Range("N58").Value = -4134237.5
$ws.Range("H59").Value = 23000
$ws.Range("J59").Value = 23000
$ws.Range("L59").Value = 23000
$ws.Range("N59").Value = -25290
$ws.Range("H60").Value = 13762.25
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 17524.5
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 17524.5
$ws.Range("N60").Value = -18546.5
$ws.Range("M60").Value = -9489
$ws.Range("H61").Value = 17433
$ws.Range("J61").Value = 17433
$ws.Range("L61").Value = 17433
$ws.Range("N61").Value = -18129
$ws.Range("H68").Value = 20500
$ws.Range("J68").Value = 20500
$ws.Range("L68").Value = 20500
$ws.Range("N68").Value = -21998
$ws.Range("H71").Value = 20500
$ws.Range("J71").Value = 20500
$ws.Range("L71").Value = 61500
$ws.Range("N71").Value = -68988
$ws.Range("H74").Value = 24264.166
$ws.Range("J74").Value = 28460
$ws.Range("L74").Value = 28460
$ws.Range("N74").Value = -30208
$ws.Range("H77").Value = 24264.166
$ws.Range("J77").Value = 28460
$ws.Range("L77").Value = 85380
$ws.Range("N77").Value = -94116
$ws.Range("H99").Value = 8698.82
$ws.Range("I99").Value = 16637.25
$ws.Range("J99").Value = 6650.1934
$ws.Range("K99").Value = 16637.25
$ws.Range("L99").Value = 6650.1934
$ws.Range("M99").Value = -15139.25
$ws.Range("N99").Value = -9646.1934
$ws.Range("H122").Value = 3406.762
$ws.Range("I122").Value = 4661.696
$ws.Range("J122").Value = 1887.6316
$ws.Range("K122").Value = 13985.088
$ws.Range("L122").Value = 5662.8948
$ws.Range("M122").Value = -11535.088
$ws.Range("N122").Value = -10562.8948
$ws.Range("H126").Value = 8698.82
$ws.Range("I126").Value = 16637.25
$ws.Range("J126").Value = 6650.1934
$ws.Range("K126").Value = 49911.75
$ws.Range("L126").Value = 19950.5802
$ws.Range("M126").Value = -47441.75
$ws.Range("N126").Value = -24890.5802
$ws.Range("H134").Value = 977800.0600000001
$ws.Range("I134").Value = 1407.3235
$ws.Range("J134").Value = 5720279
$ws.Range("K134").Value = 4221.970499999999
$ws.Range("L134").Value = 17160837
$ws.Range("M134").Value = -1686.970499999999
$ws.Range("N134").Value = -17165907
$ws.Range("H136").Value = 1757048.5
$ws.Range("I136").Value = 14074.333
$ws.Range("J136").Value = 4133831.5
$ws.Range("K136").Value = 42222.999
$ws.Range("L136").Value = 12401494.5
$ws.Range("M136").Value = -39672.999
$ws.Range("N136").Value = -12406594.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1026546.5
$ws.Range("I107").Value = 2849141
$ws.Range("J107").Value = 1337.125
$ws.Range("K107").Value = 8547423
$ws.Range("L107").Value = 4011.375
$ws.Range("M107").Value = -8545503
$ws.Range("N107").Value = -7851.375
$ws.Range("H137").Value = 8271.682000000001
$ws.Range("I137").Value = 3493.3333
$ws.Range("J137").Value = 10063.5625
$ws.Range("K137").Value = 10479.9999
$ws.Range("L137").Value = 30190.6875
$ws.Range("M137").Value = -5379.999899999999
$ws.Range("N137").Value = -40390.6875

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6312.3335
$ws.Range("I122").Value = 4940
$ws.Range("J122").Value = 8027.75
$ws.Range("K122").Value = 14820
$ws.Range("L122").Value = 24083.25
$ws.Range("M122").Value = -12370
$ws.Range("N122").Value = -28983.25
$ws.Range("H132").Value = 6910203
$ws.Range("I132").Value = 6350368.5
$ws.Range("J132").Value = 9093557
$ws.Range("K132").Value = 19051105.5
$ws.Range("L132").Value = 27280671
$ws.Range("M132").Value = -19048575.5
$ws.Range("N132").Value = -27285731

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17858638
$ws.Range("I22").Value = 833.3333
$ws.Range("J22").Value = 22728948
$ws.Range("K22").Value = 833.3333
$ws.Range("L22").Value = 22728948
$ws.Range("M22").Value = -538.3333
$ws.Range("N22").Value = -22729538
$ws.Range("H27").Value = 17858638
$ws.Range("I27").Value = 833.3333
$ws.Range("J27").Value = 22728948
$ws.Range("K27").Value = 833.3333
$ws.Range("L27").Value = 22728948
$ws.Range("M27").Value = -726.3333
$ws.Range("N27").Value = -22729162
$ws.Range("H61").Value = 4581
$ws.Range("I61").Value = 3443.6667
$ws.Range("J61").Value = 5068.4287
$ws.Range("K61").Value = 3443.6667
$ws.Range("L61").Value = 5068.4287
$ws.Range("M61").Value = -3241.6667
$ws.Range("N61").Value = -5472.4287
$ws.Range("H113").Value = 4581
$ws.Range("I113").Value = 3443.6667
$ws.Range("J113").Value = 5068.4287
$ws.Range("K113").Value = 3443.6667
$ws.Range("L113").Value = 5068.4287
$ws.Range("M113").Value = -1273.6667
$ws.Range("N113").Value = -9408.4287
$ws.Range("H132").Value = 3489292.5
$ws.Range("I132").Value = 4468453
$ws.Range("K132").Value = 13405359
$ws.Range("M132").Value = -13402829
$ws.Range("H136").Value = 5001700
$ws.Range("I136").Value = 9616255
$ws.Range("J136").Value = 2598.75
$ws.Range("K136").Value = 28848765
$ws.Range("L136").Value = 7796.25
$ws.Range("M136").Value = -28846215
$ws.Range("N136").Value = -12896.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 100
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 1356.65
$ws.Range("I122").Value = 1156.3846
$ws.Range("J122").Value = 1728.5714
$ws.Range("K122").Value = 3469.1538
$ws.Range("L122").Value = 5185.7142
$ws.Range("M122").Value = -1019.1538
$ws.Range("N122").Value = -10085.7142

